$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F306").Value = 74483
$ws.Range("G306").Value = 7358
$ws.Range("F307").Value = 73076
$ws.Range("G307").Value = 6121
$ws.Range("F309").Value = 75001
$ws.Range("G309").Value = 5311
$ws.Range("F310").Value = 76478
$ws.Range("G310").Value = 3931
$ws.Range("F311").Value = 60528
$ws.Range("G311").Value = 1887
$ws.Range("F313").Value = 74449
$ws.Range("G313").Value = 3294
$ws.Range("F314").Value = 63541
$ws.Range("G314").Value = 3029
$ws.Range("F334").Value = 192921
$ws.Range("F385").Value = 151273
$ws.Range("F386").Value = 183526
$ws.Range("F387").Value = 351949
$ws.Range("F388").Value = 731468
$ws.Range("F390").Value = 220267
$ws.Range("F391").Value = 178270
$ws.Range("F399").Value = 200830
$ws.Range("F400").Value = 150052
$ws.Range("F401").Value = 273279
$ws.Range("F413").Value = 149968
$ws.Range("F433").Value = 87465
$ws.Range("F454").Value = 52768
$ws.Range("G454").Value = 134
$ws.Range("F471").Value = 67018
$ws.Range("F472").Value = 51980
$ws.Range("F473").Value = 40227
$ws.Range("G473").Value = 42
$ws.Range("F474").Value = 45950
$ws.Range("F475").Value = 36687
$ws.Range("F476").Value = 37530
$ws.Range("F477").Value = 37047
$ws.Range("G477").Value = 36
$ws.Range("F478").Value = 55137
$ws.Range("F479").Value = 42734
$ws.Range("F480").Value = 34018
$ws.Range("F481").Value = 41502
$ws.Range("F482").Value = 36514
$ws.Range("F483").Value = 65643
$ws.Range("F485").Value = 14010
$ws.Range("F486").Value = 8938
$ws.Range("F488").Value = 6340
$ws.Range("F491").Value = 9892
$ws.Range("F492").Value = 14126
$ws.Range("F493").Value = 8296
$ws.Range("G493").Value = 9
$ws.Range("F494").Value = 6688
$ws.Range("F495").Value = 10358
$ws.Range("F496").Value = 8203
$ws.Range("F498").Value = 9102
$ws.Range("F499").Value = 11300
$ws.Range("F500").Value = 7639
$ws.Range("F501").Value = 5679
$ws.Range("F503").Value = 7381
$ws.Range("F504").Value = 7421
$ws.Range("F505").Value = 8461
$ws.Range("F506").Value = 10721
$ws.Range("F507").Value = 7132
$ws.Range("F508").Value = 5607
$ws.Range("G508").Value = 12
$ws.Range("F509").Value = 9464
$ws.Range("F510").Value = 7712
$ws.Range("F511").Value = 6646
$ws.Range("F512").Value = 8305
$ws.Range("F513").Value = 10039
$ws.Range("F514").Value = 6650
$ws.Range("G514").Value = 11
$ws.Range("F515").Value = 4840
$ws.Range("G515").Value = 14

$ws.Range("A516").Value = 44410
$ws.Range("B516").Value = 392751
$ws.Range("C516").Value = 7703
$ws.Range("D516").Value = 41
$ws.Range("E516").Value = 12541
$ws.Range("F516").Value = 8413
$ws.Range("G516").Value = 29

$ws.Range("A517").Value = 44411
$ws.Range("B517").Value = 392845
$ws.Range("C517").Value = 6918
$ws.Range("D517").Value = 94
$ws.Range("E517").Value = 12541
$ws.Range("F517").Value = 5985
$ws.Range("G517").Value = 13

$ws.Range("A518").Value = 44412
$ws.Range("B518").Value = 392898
$ws.Range("C518").Value = 5837
$ws.Range("D518").Value = 53
$ws.Range("E518").Value = 12541
$ws.Range("F518").Value = 4733
$ws.Range("G518").Value = 6

